$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"1.988074333333333"
$ws.Range("H2").Value = [double]"5.964223"
$ws.Range("I2").Value = [double]"0.01657769708907969"
$ws.Range("J2").Value = [double]"0.01657769708907968"
$ws.Range("M2").Value = [double]"0.02270466666666667"
$ws.Range("N2").Value = [double]"0.06811400000000001"
$ws.Range("O2").Value = [double]"0.002206225855740089"
$ws.Range("P2").Value = [double]"0.002206225855740089"
$ws.Range("Q2").Value = [double]"0.04513856504688889"
$ws.Range("R2").Value = [double]"0.406247085422"
$ws.Range("S2").Value = [double]"3.657414394655481E-05"
$ws.Range("T2").Value = [double]"3.65741439465548E-05"
$ws.Range("G3").Value = [double]"1.988074333333333"
$ws.Range("H3").Value = [double]"5.964223"
$ws.Range("I3").Value = [double]"0.01657769708907969"
$ws.Range("J3").Value = [double]"0.01657769708907968"
$ws.Range("O3").Value = [double]"0.002281111990432972"
$ws.Range("P3").Value = [double]"0.002281111990432972"
$ws.Range("Q3").Value = [double]"0.04667070766644445"
$ws.Range("R3").Value = [double]"0.420036368998"
$ws.Range("S3").Value = [double]"3.781558360366546E-05"
$ws.Range("T3").Value = [double]"3.781558360366545E-05"
$ws.Range("G4").Value = [double]"1.988074333333333"
$ws.Range("H4").Value = [double]"5.964223"
$ws.Range("I4").Value = [double]"0.01657769708907969"
$ws.Range("J4").Value = [double]"0.01657769708907968"
$ws.Range("M4").Value = [double]"10.24499966666667"
$ws.Range("N4").Value = [double]"30.734999"
$ws.Range("O4").Value = [double]"0.9955126621538269"
$ws.Range("P4").Value = [double]"0.9955126621538269"
$ws.Range("Q4").Value = [double]"20.36782088230855"
$ws.Range("R4").Value = [double]"183.310387940777"
$ws.Range("S4").Value = [double]"0.01650330736152947"
$ws.Range("T4").Value = [double]"0.01650330736152946"
$ws.Range("I5").Value = [double]"0.7746030815641455"
$ws.Range("J5").Value = [double]"0.7746030815641454"
$ws.Range("M5").Value = [double]"0.02270466666666667"
$ws.Range("N5").Value = [double]"0.06811400000000001"
$ws.Range("O5").Value = [double]"0.002206225855740089"
$ws.Range("P5").Value = [double]"0.002206225855740089"
$ws.Range("Q5").Value = [double]"2.109127184241778"
$ws.Range("R5").Value = [double]"18.982144658176"
$ws.Range("S5").Value = [double]"0.001708949346482766"
$ws.Range("T5").Value = [double]"0.001708949346482766"
$ws.Range("I6").Value = [double]"0.7746030815641455"
$ws.Range("J6").Value = [double]"0.7746030815641454"
$ws.Range("O6").Value = [double]"0.002281111990432972"
$ws.Range("P6").Value = [double]"0.002281111990432972"
$ws.Range("S6").Value = [double]"0.001766956377182302"
$ws.Range("T6").Value = [double]"0.001766956377182302"
$ws.Range("I7").Value = [double]"0.7746030815641455"
$ws.Range("J7").Value = [double]"0.7746030815641454"
$ws.Range("M7").Value = [double]"10.24499966666667"
$ws.Range("N7").Value = [double]"30.734999"
$ws.Range("O7").Value = [double]"0.9955126621538269"
$ws.Range("P7").Value = [double]"0.9955126621538269"
$ws.Range("Q7").Value = [double]"951.6989443953352"
$ws.Range("R7").Value = [double]"8565.290499558016"
$ws.Range("S7").Value = [double]"0.7711271758404804"
$ws.Range("T7").Value = [double]"0.7711271758404803"
$ws.Range("G8").Value = [double]"23.741365"
$ws.Range("H8").Value = [double]"71.22409500000001"
$ws.Range("I8").Value = [double]"0.1979690350870239"
$ws.Range("J8").Value = [double]"0.1979690350870239"
$ws.Range("M8").Value = [double]"0.02270466666666667"
$ws.Range("N8").Value = [double]"0.06811400000000001"
$ws.Range("O8").Value = [double]"0.002206225855740089"
$ws.Range("P8").Value = [double]"0.002206225855740089"
$ws.Range("Q8").Value = [double]"0.5390397785366667"
$ws.Range("R8").Value = [double]"4.851358006830001"
$ws.Range("S8").Value = [double]"0.000436764403844909"
$ws.Range("T8").Value = [double]"0.000436764403844909"
$ws.Range("G9").Value = [double]"23.741365"
$ws.Range("H9").Value = [double]"71.22409500000001"
$ws.Range("I9").Value = [double]"0.1979690350870239"
$ws.Range("J9").Value = [double]"0.1979690350870239"
$ws.Range("O9").Value = [double]"0.002281111990432972"
$ws.Range("P9").Value = [double]"0.002281111990432972"
$ws.Range("Q9").Value = [double]"0.5573364571633334"
$ws.Range("R9").Value = [double]"5.016028114470001"
$ws.Range("S9").Value = [double]"0.0004515895396714561"
$ws.Range("T9").Value = [double]"0.0004515895396714561"
$ws.Range("G10").Value = [double]"23.741365"
$ws.Range("H10").Value = [double]"71.22409500000001"
$ws.Range("I10").Value = [double]"0.1979690350870239"
$ws.Range("J10").Value = [double]"0.1979690350870239"
$ws.Range("M10").Value = [double]"10.24499966666667"
$ws.Range("N10").Value = [double]"30.734999"
$ws.Range("O10").Value = [double]"0.9955126621538269"
$ws.Range("P10").Value = [double]"0.9955126621538269"
$ws.Range("Q10").Value = [double]"243.2302765112117"
$ws.Range("R10").Value = [double]"2189.072488600905"
$ws.Range("S10").Value = [double]"0.1970806811435076"
$ws.Range("T10").Value = [double]"0.1970806811435076"
$ws.Range("G11").Value = [double]"1.301204666666667"
$ws.Range("H11").Value = [double]"3.903614"
$ws.Range("I11").Value = [double]"0.01085018625975097"
$ws.Range("J11").Value = [double]"0.01085018625975097"
$ws.Range("M11").Value = [double]"0.02270466666666667"
$ws.Range("N11").Value = [double]"0.06811400000000001"
$ws.Range("O11").Value = [double]"0.002206225855740089"
$ws.Range("P11").Value = [double]"0.002206225855740089"
$ws.Range("Q11").Value = [double]"0.02954341822177778"
$ws.Range("R11").Value = [double]"0.265890763996"
$ws.Range("S11").Value = [double]"2.393796146585844E-05"
$ws.Range("T11").Value = [double]"2.393796146585843E-05"
$ws.Range("G12").Value = [double]"1.301204666666667"
$ws.Range("H12").Value = [double]"3.903614"
$ws.Range("I12").Value = [double]"0.01085018625975097"
$ws.Range("J12").Value = [double]"0.01085018625975097"
$ws.Range("O12").Value = [double]"0.002281111990432972"
$ws.Range("P12").Value = [double]"0.002281111990432972"
$ws.Range("Q12").Value = [double]"0.03054621328488889"
$ws.Range("R12").Value = [double]"0.274915919564"
$ws.Range("S12").Value = [double]"2.475048997554903E-05"
$ws.Range("T12").Value = [double]"2.475048997554902E-05"
$ws.Range("G13").Value = [double]"1.301204666666667"
$ws.Range("H13").Value = [double]"3.903614"
$ws.Range("I13").Value = [double]"0.01085018625975097"
$ws.Range("J13").Value = [double]"0.01085018625975097"
$ws.Range("M13").Value = [double]"10.24499966666667"
$ws.Range("N13").Value = [double]"30.734999"
$ws.Range("O13").Value = [double]"0.9955126621538269"
$ws.Range("P13").Value = [double]"0.9955126621538269"
$ws.Range("Q13").Value = [double]"13.33084137626511"
$ws.Range("R13").Value = [double]"119.977572386386"
$ws.Range("S13").Value = [double]"0.01080149780830957"
$ws.Range("T13").Value = [double]"0.01080149780830956"
